$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edited cell is C13 (whose value is mirrored across C13:F13).
# Select it first so the saved sheet view reflects the user's edit location.
$ws.Range("C13").Select() | Out-Null

# Previously this row held "vdDACg67q2GZB0XbKQc"; update it to
# "HvdDACg67q2GZB0XbKQc" across the merged/repeated values in C13:F13.
$ws.Range("C13:F13").Value = "HvdDACg67q2GZB0XbKQc"
